$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, `
                               $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find.Execute failed for: $find"
    }
}

# 1. Precondiciones: append sentence about owner being registered, drop the
#    trailing lone space run.
Replace-Text "dentro de este. " "dentro de este. El propietario está registrado en el sistema."

# 2. Step: "El sistema muestra la pantalla correspondiente al registro de
#    pacientes." now appears earlier (step 3) -- do this BEFORE touching step 1
#    text, to avoid collisions with the new text we introduce in step 1 below.
Replace-Text "El sistema muestra la pantalla correspondiente al registro de pacientes." "El sistema muestra los datos del propietario para tener una confirmación visual sobre la existencia del propietario en el sistema."

# 3. Step 1 of flujo básico: selecting menu option -> system shows screen.
Replace-Text "El usuario selecciona en el menú principal la opción “Gestión de pacientes”." "El sistema muestra la pantalla correspondiente al registro de pacientes."

# 4. Step 2: "En el menú nuevo..." -> "El usuario ingresa el documento del propietario."
Replace-Text "En el menú nuevo selecciona la opción “Registrar paciente”." "El usuario ingresa el documento del propietario."

# 5. Step 4: remove "que son el documento del propietario" clause.
Replace-Text " ingresa los datos básicos del paciente que son el documento del propietario, el nombre" " ingresa los datos básicos del paciente como lo son el nombre"

# 6. Step 5: system verification no longer checks owner document existence.
Replace-Text "El sistema verifica que el documento del propietario corresponda a un propietario registrado, que el nombre" "El sistema verifica que el nombre"

# 7. Remove the "El usuario es redirigido al menú principal." step paragraph entirely.
$removed = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "El usuario es redirigido al menú principal.") {
        $p.Range.Delete()
        $removed = $true
        break
    }
}
if (-not $removed) {
    throw "Could not find paragraph to delete: El usuario es redirigido al menú principal."
}

# 8. Post-condiciones: add "correctamente".
Replace-Text "Se ha registrado la información del paciente." "Se ha registrado la información del paciente correctamente."

# 9. Excepción "Propietario no registrado": step number 4 -> 3.
Replace-Text "En el paso 4 del flujo normal : si documento del propietario no existe" "En el paso 3 del flujo normal : si documento del propietario no existe"

# 10. Excepción "Campo no diligenciado": step number 4 -> 5 (note no space
#     before the digit in the original text).
Replace-Text "En el paso4 del flujo normal de eventos" "En el paso 5 del flujo normal de eventos"
